$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "59.441.84"
$ws.Range("E2").Value = "  +0.45%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.525.39"
$ws.Range("E3").Value = "  +0.60%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "536.97"
$ws.Range("E5").Value = "  -0.38%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "139.86"
$ws.Range("E6").Value = "  -3.32%  "
$ws.Range("E7").Value = "  +0.36%  "
$ws.Range("E8").Value = "  -1.72%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.531.56"
$ws.Range("E9").Value = "  -0.57%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0994"
$ws.Range("E10").Value = "  -0.59%  "
$ws.Range("E11").Value = "  +1.22%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.40"
$ws.Range("E12").Value = "  -2.83%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.355"
$ws.Range("E13").Value = "  +0.73%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.972.75"
$ws.Range("E14").Value = "  +0.65%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "23.14"
$ws.Range("E15").Value = "  -2.59%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "59.339.73"
$ws.Range("E16").Value = "  +0.46%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000141"
$ws.Range("E17").Value = "  +0.71%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.530.96"
$ws.Range("E18").Value = "  -0.28%  "
$ws.Range("E19").Value = "  -3.47%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.22"
$ws.Range("E20").Value = "  -1.63%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "321.82"
$ws.Range("E21").Value = "  -0.61%  "
$ws.Range("E22").Value = "  +0.11%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.80"
$ws.Range("E23").Value = "  +0.27%  "
$ws.Range("E24").Value = "  -1.25%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.422"
$ws.Range("E25").Value = "  -4.44%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.166"
$ws.Range("E26").Value = "  +1.71%  "
$ws.Range("E27").Value = "  +0.72%  "
$ws.Range("E28").Value = "  +0.05%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.75"
$ws.Range("E29").Value = "  -0.49%  "
$ws.Range("B30").Value = "PancakeSwap"
$ws.Range("C30").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.80"
$ws.Range("E30").Value = "  +0.44%  "
$ws.Range("B31").Value = "PEPE"
$ws.Range("C31").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0₃0767"
$ws.Range("E31").Value = "  -1.53%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "160.90"
$ws.Range("E32").Value = "  +1.47%  "
$ws.Range("E34").Value = "  +0.76%  "
$ws.Range("E35").Value = "  -7.58%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "18.55"
$ws.Range("E36").Value = "  -0.32%  "
$ws.Range("E37").Value = "  -4.79%  "
$ws.Range("E38").Value = "  -1.98%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "37.01"
$ws.Range("E39").Value = "  +0.26%  "
$ws.Range("E40").Value = "  -0.50%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.809"
$ws.Range("E41").Value = "  -1.94%  "
$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.27"
$ws.Range("E42").Value = "  -8.53%  "
$ws.Range("B43").Value = "Bittensor"
$ws.Range("C43").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "285.02"
$ws.Range("E43").Value = "  -6.62%  "
$ws.Range("E44").Value = "  +0.44%  "
$ws.Range("E45").Value = "  -0.75%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.85"
$ws.Range("E46").Value = "  +0.64%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "124.04"
$ws.Range("E47").Value = "  -1.83%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "18.55"
$ws.Range("E49").Value = "  -1.04%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0510"
$ws.Range("E50").Value = "  -1.71%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0223"
$ws.Range("E51").Value = "  -2.40%  "
